# ---------------------------------------------------------------------
# top15_perfecte_woningen_tabel_final.xlsx - rebuild table with new
# "Tuin" column and refreshed scoring algorithm output.
# ---------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old "Score" header one column to the right (G1 -> H1),
# re-using the existing header style (bold font + thin border +
# centered alignment) by copying the formatting, then put the new
# "Tuin" header into the now-vacated G1 using the same style.
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)
$ws.Cells.Item(1, 8).Value = "Score"
$ws.Cells.Item(1, 7).Value = "Tuin"

# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Nieuwe Leliestraat 97, 1015SL, Amsterdam"
$ws.Cells.Item(2, 3).Value = 925000
$ws.Cells.Item(2, 4).Value = 100
$ws.Cells.Item(2, 5).Value = "Unknown"
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(2, 7).Value = ""
$ws.Cells.Item(2, 8).Value = 0.8392333333333335

# Row 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Nieuwe Leliestraat 8H, 1015SP, Amsterdam"
$ws.Cells.Item(3, 3).Value = 895000
$ws.Cells.Item(3, 4).Value = 98
$ws.Cells.Item(3, 5).Value = "Unknown"
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = ""
$ws.Cells.Item(3, 8).Value = 0.8192333333333334

# Row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Oude Looiersstraat 382, 1016VJ, Amsterdam"
$ws.Cells.Item(4, 3).Value = 795000
$ws.Cells.Item(4, 4).Value = 90
$ws.Cells.Item(4, 5).Value = "Unknown"
$ws.Cells.Item(4, 6).Value = 2
$ws.Cells.Item(4, 7).Value = ""
$ws.Cells.Item(4, 8).Value = 0.8102627450980394

# Row 5
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Tweede Bloemdwarsstraat 34A, 1016LM, Amsterdam"
$ws.Cells.Item(5, 3).Value = 750000
$ws.Cells.Item(5, 4).Value = 83
$ws.Cells.Item(5, 5).Value = "Unknown"
$ws.Cells.Item(5, 6).Value = 3
$ws.Cells.Item(5, 7).Value = ""
$ws.Cells.Item(5, 8).Value = 0.7978000000000002

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Akoleienstraat 81, 1016LN, Amsterdam"
$ws.Cells.Item(6, 3).Value = 725000
$ws.Cells.Item(6, 4).Value = 105
$ws.Cells.Item(6, 5).Value = "Unknown"
$ws.Cells.Item(6, 6).Value = 3
$ws.Cells.Item(6, 7).Value = ""
$ws.Cells.Item(6, 8).Value = 0.7966000000000001

# Row 7
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Laurierstraat 26H, 1016PM, Amsterdam"
$ws.Cells.Item(7, 3).Value = 949000
$ws.Cells.Item(7, 4).Value = 96
$ws.Cells.Item(7, 5).Value = "Unknown"
$ws.Cells.Item(7, 6).Value = 2
$ws.Cells.Item(7, 7).Value = ""
$ws.Cells.Item(7, 8).Value = 0.7940627450980393

# Row 8
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Laurierstraat 11B, 1016PG, Amsterdam"
$ws.Cells.Item(8, 3).Value = 1100000
$ws.Cells.Item(8, 4).Value = 96
$ws.Cells.Item(8, 5).Value = "Unknown"
$ws.Cells.Item(8, 6).Value = 2
$ws.Cells.Item(8, 7).Value = ""
$ws.Cells.Item(8, 8).Value = 0.7940627450980393

# Row 9
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Oude Looiersstraat 372, 1016VG, Amsterdam"
$ws.Cells.Item(9, 3).Value = 660000
$ws.Cells.Item(9, 4).Value = 83
$ws.Cells.Item(9, 5).Value = "Unknown"
$ws.Cells.Item(9, 6).Value = 2
$ws.Cells.Item(9, 7).Value = ""
$ws.Cells.Item(9, 8).Value = 0.7927627450980393

# Row 10
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Elandsgracht 103A, 1016TS, Amsterdam"
$ws.Cells.Item(10, 3).Value = 950000
$ws.Cells.Item(10, 4).Value = 113
$ws.Cells.Item(10, 5).Value = "Unknown"
$ws.Cells.Item(10, 6).Value = 3
$ws.Cells.Item(10, 7).Value = ""
$ws.Cells.Item(10, 8).Value = 0.7876294117647061

# Row 11
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Tweede Bloemdwarsstraat 92, 1016LL, Amsterdam"
$ws.Cells.Item(11, 3).Value = 799000
$ws.Cells.Item(11, 4).Value = 82
$ws.Cells.Item(11, 5).Value = "Unknown"
$ws.Cells.Item(11, 6).Value = 2
$ws.Cells.Item(11, 7).Value = ""
$ws.Cells.Item(11, 8).Value = 0.7786333333333335

# Row 12
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Nieuwe Leliestraat 168H, 1015HE, Amsterdam"
$ws.Cells.Item(12, 3).Value = 745000
$ws.Cells.Item(12, 4).Value = 126
$ws.Cells.Item(12, 5).Value = "Unknown"
$ws.Cells.Item(12, 6).Value = 3
$ws.Cells.Item(12, 7).Value = ""
$ws.Cells.Item(12, 8).Value = 0.7759000000000001

# Row 13
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "Tweede Leliedwarsstraat 16, 1015TC, Amsterdam"
$ws.Cells.Item(13, 3).Value = 775000
$ws.Cells.Item(13, 4).Value = 79
$ws.Cells.Item(13, 5).Value = "Unknown"
$ws.Cells.Item(13, 6).Value = 2
$ws.Cells.Item(13, 7).Value = ""
$ws.Cells.Item(13, 8).Value = 0.7717333333333335

# Row 14
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "Bloemstraat 1211, 1016KZ, Amsterdam"
$ws.Cells.Item(14, 3).Value = 950000
$ws.Cells.Item(14, 4).Value = 97
$ws.Cells.Item(14, 5).Value = "Unknown"
$ws.Cells.Item(14, 6).Value = 2
$ws.Cells.Item(14, 7).Value = ""
$ws.Cells.Item(14, 8).Value = 0.7693333333333334

# Row 15
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Egelantiersgracht 101A, 1015RG, Amsterdam"
$ws.Cells.Item(15, 3).Value = 850000
$ws.Cells.Item(15, 4).Value = 89
$ws.Cells.Item(15, 5).Value = "Unknown"
$ws.Cells.Item(15, 6).Value = 3
$ws.Cells.Item(15, 7).Value = ""
$ws.Cells.Item(15, 8).Value = 0.7666000000000002

# Row 16
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "Lauriergracht 1142, 1016RR, Amsterdam"
$ws.Cells.Item(16, 3).Value = 1250000
$ws.Cells.Item(16, 4).Value = 109
$ws.Cells.Item(16, 5).Value = "Unknown"
$ws.Cells.Item(16, 6).Value = 3
$ws.Cells.Item(16, 7).Value = ""
$ws.Cells.Item(16, 8).Value = 0.766429411764706
